# Append rows 17-21 of new training run results to the logs worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preprocess = 'trim "space" and ",", space after punctuation, remove break line, convert unicode to ascii, remove multiple spaces, convert to lower'
$features   = '12 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, #max_digit_skip_0 >= 7, #max_digit_skip_0 = 0, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit'
$model      = 'Neuron Network'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000'
$templateFilter = '0 filters: '
$valAccuracy = 0.940594059405941

$rows = @(
    @{ Row=17; Time='20160415_171551'; RunningTime=1844.938; TestAccuracy=0.988666666666667; TemplateFilterVal=0.146341463414634 },
    @{ Row=18; Time='20160415_174636'; RunningTime=1785.976; TestAccuracy=0.992666666666667; TemplateFilterVal=0.170731707317073 },
    @{ Row=19; Time='20160415_181622'; RunningTime=1733.836; TestAccuracy=0.993333333333333; TemplateFilterVal=0.170731707317073 },
    @{ Row=20; Time='20160415_184516'; RunningTime=1314.278; TestAccuracy=0.988666666666667; TemplateFilterVal=0.158536585365854 },
    @{ Row=21; Time='20160415_190710'; RunningTime=1204.194; TestAccuracy=0.992;              TemplateFilterVal=0.146341463414634 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Time
    $ws.Cells.Item($row, 2).Value = $r.RunningTime
    $ws.Cells.Item($row, 3).Value = $preprocess
    $ws.Cells.Item($row, 4).Value = $features
    $ws.Cells.Item($row, 5).Value = $model
    $ws.Cells.Item($row, 6).Value = $modelDetails
    $ws.Cells.Item($row, 7).Value = $r.TestAccuracy
    $ws.Cells.Item($row, 8).Value = $valAccuracy
    $ws.Cells.Item($row, 9).Value = $templateFilter
    $ws.Cells.Item($row, 10).Value = $r.TemplateFilterVal
}
